$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in C23/G23 - new commit entry with hours
$ws.Range("C23").Value = "merchant talk action update (40%)"
$ws.Range("C23").Style = $ws.Range("C17").Style
$ws.Range("G23").Value = 1.5

# Fill in C24/G24 - new commit entry with hours
$ws.Range("C24").Value = "user look_it & sell/buy update"
$ws.Range("C24").Style = $ws.Range("C17").Style
$ws.Range("G24").Value = 0.7

# Move Total row down to row 29 (capture row-26 styles first, then clear row 26)
$ws.Range("F29").Style = $ws.Range("F26").Style
$ws.Range("G29").Style = $ws.Range("G26").Style
$ws.Range("F29").Value = "Total(h):"
$ws.Range("G29").Formula = "=SUM(G4:G24)"

# Remove the old Total row content at row 26 (clear contents+format so the
# now-empty row element is dropped entirely)
$ws.Range("F26:G26").Clear()

# Update the view: move the selection (also clears the old fixed topLeftCell)
$ws.Range("G28").Select()
